$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells that contain numeric-looking text stay as text
$priceRange = $ws.Range("D2:D50")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.166.42"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").Value = "3.278.33"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "588.05"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("D6").Value = "187.00"
$ws.Range("E6").Value = "  +4.55%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("E9").Value = "  +4.31%  "

$ws.Range("D10").Value = "6.74"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("D12").Value = "3.843.67"
$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("D14").Value = "28.72"
$ws.Range("E14").Value = "  +1.69%  "

$ws.Range("D15").Value = "68.179.80"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("E16").Value = "  +2.55%  "

$ws.Range("D17").Value = "3.269.59"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("E19").Value = "  +1.73%  "

$ws.Range("D20").Value = "381.83"
$ws.Range("E20").Value = "  +1.97%  "

$ws.Range("D21").Value = "7.75"
$ws.Range("E21").Value = "  +1.46%  "

$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "71.58"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "0.0000121"
$ws.Range("E24").Value = "  +1.91%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "0.515"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("D26").Value = "9.81"
$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("E27").Value = "  +5.15%  "

$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").Value = "5.83"
$ws.Range("E29").Value = "  +3.46%  "

$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "22.96"
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("D32").Value = "7.19"
$ws.Range("E32").Value = "  +5.41%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.28"
$ws.Range("E33").Value = "  +1.21%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +2.65%  "

$ws.Range("D36").Value = "163.08"
$ws.Range("E36").Value = "  -1.46%  "

$ws.Range("D37").Value = "1.87"
$ws.Range("E37").Value = "  -0.51%  "

$ws.Range("D38").Value = "0.840"
$ws.Range("E38").Value = "  -1.94%  "

$ws.Range("E39").Value = "  +3.69%  "

$ws.Range("D40").Value = "26.72"
$ws.Range("E40").Value = "  -1.13%  "

$ws.Range("D41").Value = "4.61"
$ws.Range("E41").Value = "  +4.73%  "

$ws.Range("D42").Value = "2.64"
$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "25.56"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "347.64"
$ws.Range("E44").Value = "  -1.97%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0691"
$ws.Range("E45").Value = "  +2.26%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "41.30"
$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("D47").Value = "2.650.48"
$ws.Range("E47").Value = "  -4.18%  "

$ws.Range("E48").Value = "  +1.91%  "

$ws.Range("D49").Value = "32.30"
$ws.Range("E49").Value = "  +4.19%  "

$ws.Range("E50").Value = "  +1.50%  "

# Restore normal (unstyled) formatting for the price column so only values changed
$priceRange.Style = "Normal"

Write-Host "Done applying crypto list updates"